$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set D and E columns to Text format before assigning, to preserve values
# exactly as strings (Excel would otherwise auto-convert numeric-looking
# text into real numbers and strip formatting like leading zeros / trailing dots).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.533.97"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.81%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.569.50"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.35%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -1.59%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.92"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.10%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.61%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.67"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.70%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.251"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.41%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0596"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.27%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.23%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.790.35"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.25%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.567.86"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.06%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.68%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.04%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.533.63"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.89%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.48"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.93%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "225.15"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +4.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.52"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.77%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0704"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.990"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.51%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.13"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.43"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.28%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.32%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.80"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.60%  "
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.20"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.78%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.62"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.20%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.108"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.57%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0471"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.65%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.28%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.451.37"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.88%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.16"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.44%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.40%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.40%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.99%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0168"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.543"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.93%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.73"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.05%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.81%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.990"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +5.37%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.973"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.95%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "64.68"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.703.49"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.64"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.02%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.49%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₇0993"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -3.38%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0948"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.32%  "
